# Update header labels across sheets so the first row can be used as a
# header automatically when imported into Power BI.

$wb = $excel.ActiveWorkbook

# Sheets that use the "Ano" prefix pattern (year columns B:E)
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Sheet that uses "Intervalo" prefix pattern (interval columns B:E)
$wsIncremental = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIncremental.Range("B1").Value = "Intervalo 2015"
$wsIncremental.Range("C1").Value = "Intervalo 2015-2030"
$wsIncremental.Range("D1").Value = "Intervalo 2031-2040"
$wsIncremental.Range("E1").Value = "Intervalo 2041-2050"

# Sheet with only a single year column (B1)
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Range("B1").Value = "Ano 2015"
